# Format the transferred_at column (A) as real dates instead of
# text-based shared-string values, so the sheet can be used for
# import/export round-tripping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and row 3 both hold the same transferred_at date: 2020-01-01,
# which as an Excel serial date number is 43831. Replace the text
# values with real numeric date values and apply a yyyy-mm-dd date
# number format to them.
$dateRange = $ws.Range("A2:A3")
$dateRange.Value = 43831
$dateRange.NumberFormat = "yyyy\-mm\-dd"

# Restore the sheet's active cell/selection like the original commit.
$ws.Range("D25").Select() | Out-Null
